# Add team record (Wins / Losses / Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of the last existing header cell (AC1, style index 1:
# bold, centered/top-aligned, thin border) onto the three new header cells
# so AD1:AF1 look consistent with the rest of the header row.
$ws.Range("AC1").Copy($ws.Range("AD1:AF1"))

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-45) gets the same team record values.
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 30).Value = 100  # AD: Wins
    $ws.Cells.Item($r, 31).Value = 62   # AE: Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF: Ties
}
